# Excel File reader POC
# - Adds a new "final" worksheet at the end of the workbook, containing a
#   combined header row (SR.NO / ESIC User / ESIC Password / ... employee
#   columns), formatted with a bold/centered style.
# - Adds two new shared strings ("ESIC User", "ESIC Password").
# - Updates the selection on Sheet1 (no longer the tab shown when the file
#   is opened) and makes the new "final" sheet the active / selected tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update Sheet1's selection (it is no longer the active tab once the
#    new sheet is added, so select it first while it is still active).
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("C3").Select()

# ---------------------------------------------------------------------
# 2. Add the new "final" worksheet at the end of the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "final"

# ---------------------------------------------------------------------
# 3. Populate the header row (row 1) with the combined set of columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "SR. NO."
$ws.Range("B1").Value = "ESIC User"
$ws.Range("C1").Value = "ESIC Password"
$ws.Range("D1").Value = "EMP CODE"
$ws.Range("E1").Value = "INS NO"
$ws.Range("F1").Value = "EMPLOYEE NAME"
$ws.Range("G1").Value = "Gender"
$ws.Range("H1").Value = "DATE OF BIRTH"
$ws.Range("I1").Value = "FATHER NAME"
$ws.Range("J1").Value = "DATE OF BIRTH OF FATHER.HUSBAND"
$ws.Range("K1").Value = "NAME OF MOTHER"
$ws.Range("L1").Value = "DATE OF BIRTH OF MOTHER"
$ws.Range("M1").Value = "MARITAL STATUS"
$ws.Range("N1").Value = "NAME OF SPOUSE"
$ws.Range("O1").Value = "DATE OF BIRTH OF SPOUSE"
$ws.Range("P1").Value = "NAME OF SON"
$ws.Range("Q1").Value = "DATE OF BIRTH OF SON"
$ws.Range("R1").Value = "NAME OF DAUGHTER "
$ws.Range("S1").Value = "DATE OF BIRTH OF DAUGHTER"
$ws.Range("T1").Value = "LOCAL ADDRESS"
$ws.Range("U1").Value = "PIN CODE"
$ws.Range("V1").Value = "PERMANENT ADDRESS"
$ws.Range("W1").Value = "PIN CODE"
$ws.Range("X1").Value = "CONTACT NO"
$ws.Range("Y1").Value = "DATE OF APPOINTMENT"
$ws.Range("Z1").Value = "NAME OF NOMINEE "
$ws.Range("AA1").Value = "RELATIONSHIP WITH THE EMPLOYEE"

# ---------------------------------------------------------------------
# 4. Format the header row: bold font, centered horizontally/vertically.
#    Format A1 first, then copy/paste the resulting format onto the rest
#    of the row so every cell shares a single new style entry.
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1").VerticalAlignment = -4108     # xlCenter

$ws.Range("A1").Copy()
$ws.Range("B1:AA1").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(1).RowHeight = 42

# ---------------------------------------------------------------------
# 5. Approximate the column widths from the original workbook (best
#    effort; the exact fractional widths cannot always be reproduced).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 6.74
$ws.Columns.Item(3).ColumnWidth = 12.88
$ws.Columns.Item(4).ColumnWidth = 9.31
$ws.Columns.Item(5).ColumnWidth = 6.45
$ws.Columns.Item(6).ColumnWidth = 15.31
$ws.Columns.Item(7).ColumnWidth = 6.74
$ws.Columns.Item(8).ColumnWidth = 13.31
$ws.Columns.Item(9).ColumnWidth = 12.88
$ws.Columns.Item(10).ColumnWidth = 33.31
$ws.Columns.Item(11).ColumnWidth = 16.74
$ws.Columns.Item(12).ColumnWidth = 24.45
$ws.Columns.Item(13).ColumnWidth = 15.17
$ws.Columns.Item(14).ColumnWidth = 15.88
$ws.Columns.Item(15).ColumnWidth = 23.74
$ws.Columns.Item(16).ColumnWidth = 12.88
$ws.Columns.Item(17).ColumnWidth = 20.74
$ws.Columns.Item(18).ColumnWidth = 19.17
$ws.Columns.Item(19).ColumnWidth = 26.59
$ws.Columns.Item(20).ColumnWidth = 14.31
$ws.Columns.Item(21).ColumnWidth = 8.59
$ws.Columns.Item(22).ColumnWidth = 20.02
$ws.Columns.Item(23).ColumnWidth = 8.59
$ws.Columns.Item(24).ColumnWidth = 11.88
$ws.Columns.Item(25).ColumnWidth = 21.59
$ws.Columns.Item(26).ColumnWidth = 18.17
$ws.Columns.Item(27).ColumnWidth = 32.02

# ---------------------------------------------------------------------
# 6. Page setup for the new sheet.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

# ---------------------------------------------------------------------
# 7. Make the new sheet the active / selected sheet (matches the
#    workbook's new activeTab pointing at "final").
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
